$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value of 45189 (2023-09-20)
# for every data row (2 through 351). Update it to 45190 (2023-09-21).
$startRow = 2
$endRow = 351

for ($r = $startRow; $r -le $endRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45190
}
